$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '42.355.42'
$ws.Range("E2").Value = '  +1.56%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.177.04'
$ws.Range("E3").Value = '  +0.19%  '

$ws.Range("E4").Value = '  -0.14%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '253.24'
$ws.Range("E5").Value = '  +6.47%  '

$ws.Range("E6").Value = '  -0.09%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '75.16'
$ws.Range("E7").Value = '  +4.11%  '

$ws.Range("E8").Value = '  -0.02%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.583'
$ws.Range("E9").Value = '  +0.50%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '41.09'
$ws.Range("E10").Value = '  +3.22%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0912'
$ws.Range("E11").Value = '  +0.47%  '

$ws.Range("B12").Value = 'TRON'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.101'
$ws.Range("E12").Value = '  +0.98%  '

$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.77'
$ws.Range("E13").Value = '  +1.11%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.504.75'
$ws.Range("E14").Value = '  +0.02%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '14.17'
$ws.Range("E15").Value = '  -1.34%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.183.66'
$ws.Range("E16").Value = '  +0.17%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.768'
$ws.Range("E17").Value = '  -1.29%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '42.271.12'
$ws.Range("E18").Value = '  +1.70%  '

$ws.Range("E19").Value = '  -0.27%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '70.67'
$ws.Range("E20").Value = '  +0.93%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.87'
$ws.Range("E21").Value = '  +1.42%  '

$ws.Range("B22").Value = 'ImmutableX'
$ws.Range("C22").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.18'
$ws.Range("E22").Value = '  +6.06%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '226.93'
$ws.Range("E23").Value = '  +0.49%  '

$ws.Range("B24").Value = 'InternetComputer(DFINITY)'
$ws.Range("C24").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.57'
$ws.Range("E24").Value = '  -4.13%  '

$ws.Range("E25").Value = '  -0.13%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '10.45'
$ws.Range("E26").Value = '  -2.08%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.35'
$ws.Range("E27").Value = '  +2.60%  '

$ws.Range("B28").Value = 'Toncoin'
$ws.Range("C28").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.22'
$ws.Range("E28").Value = '  +6.75%  '

$ws.Range("B29").Value = 'PancakeSwap'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.18'
$ws.Range("E29").Value = '  -0.21%  '

$ws.Range("B30").Value = 'InjectiveProtocol'
$ws.Range("C30").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '37.08'
$ws.Range("E30").Value = '  +12.39%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '169.32'
$ws.Range("E31").Value = '  -1.33%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '20.00'
$ws.Range("E32").Value = '  +0.80%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0816'
$ws.Range("E33").Value = '  +5.80%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.11'
$ws.Range("E34").Value = '  -3.06%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.120'
$ws.Range("E35").Value = '  +0.06%  '

$ws.Range("E36").Value = '  +4.59%  '

$ws.Range("E37").Value = '  -0.69%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0334'
$ws.Range("E38").Value = '  +8.47%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '11.85'
$ws.Range("E39").Value = '  -0.97%  '

$ws.Range("E40").Value = '  -1.32%  '

$ws.Range("E41").Value = '  +4.14%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '59.54'
$ws.Range("E42").Value = '  +1.00%  '

$ws.Range("E43").Value = '  -3.67%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '103.02'
$ws.Range("E44").Value = '  +6.15%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.472'
$ws.Range("E45").Value = '  +15.59%  '

$ws.Range("E46").Value = '  +11.09%  '

$ws.Range("E47").Value = '  -1.71%  '

$ws.Range("E48").Value = '  +0.71%  '

$ws.Range("E49").Value = '  +1.00%  '

$ws.Range("E50").Value = '  +1.55%  '

$ws.Range("E51").Value = '  +0.48%  '
